$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right after "2021-Q4", before the
#    existing "总计" (totals) sheet.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Fetch the totals sheet only *after* inserting the new sheet above — the
# insertion shifts everyone after "2021-Q4" down by one position, so a
# reference grabbed beforehand would (incorrectly) keep pointing at the
# freshly-added "2022-Q1" sheet's slot instead of following "总计" sheet.
$totals = $wb.Worksheets.Item("总计")

# Copy the shared header / row formatting from the already-existing
# "2021-Q4" sheet so the new sheet matches the look of its siblings.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row. The numeric-looking columns (fund code / size / position /
# weight / holding value) are stored as plain text in the source workbook,
# so force a text number format before writing them — this avoids Excel
# re-interpreting "014887" as 14887 or rounding "17.22" — and then reset
# the format back to Normal so no stray style is left on the cells.
$newSheet.Range("A2").Value = 0
$dataRow = $newSheet.Range("B2:G2")
$dataRow.NumberFormat = "@"
$newSheet.Range("B2").Value = "014887"
$newSheet.Range("C2").Value = "招商安福1年定期开放债券"
$newSheet.Range("D2").Value = "17.22"
$newSheet.Range("E2").Value = "27.65"
$newSheet.Range("F2").Value = "0.62"
$newSheet.Range("G2").Value = "0.1068"
$dataRow.Style = "Normal"
$newSheet.Range("H2").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: a new leading data row is added for
#    2022-Q1, and the pre-existing 2021-Q4 / 2021-Q3 rows move down one
#    row each, with column A renumbered (0, 1, 2).
# ---------------------------------------------------------------------------

# Give the brand-new row 4 the same look (bold/border style on column A)
# as its neighbours before filling in its values.
$totals.Range("A3").Copy()
$totals.Range("A4").PasteSpecial(-4122)

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2021-Q3"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 0.04

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2021-Q4"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.14

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.11
